$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells keep their literal string representation
# (e.g. "5.20", "91.50") instead of being auto-coerced to numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.205.67'
$ws.Range('E2').Value = '  -4.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.961.26'
$ws.Range('E3').Value = '  -4.72%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.65'
$ws.Range('E5').Value = '  -4.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.622'
$ws.Range('E6').Value = '  -4.97%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.35'
$ws.Range('E7').Value = '  -10.91%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.373'
$ws.Range('E9').Value = '  -1.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '56.86'
$ws.Range('E10').Value = '  -6.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0793'
$ws.Range('E11').Value = '  +3.66%  '
$ws.Range('E12').Value = '  -1.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.847'
$ws.Range('E13').Value = '  -8.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '13.96'
$ws.Range('E14').Value = '  -6.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.82'
$ws.Range('E15').Value = '  +5.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.247.83'
$ws.Range('E16').Value = '  -4.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.38'
$ws.Range('E17').Value = '  -3.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.962.80'
$ws.Range('E18').Value = '  -4.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '36.083.58'
$ws.Range('E19').Value = '  -4.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.12'
$ws.Range('E20').Value = '  -4.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0850'
$ws.Range('E21').Value = '  -3.51%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '236.38'
$ws.Range('E22').Value = '  -1.42%  '
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.20'
$ws.Range('E23').Value = '  -3.30%  '
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('E25').Value = '  -5.50%  '
$ws.Range('E26').Value = '  -4.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.63'
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.71'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.08'
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.125'
$ws.Range('E30').Value = '  +7.71%  '
$ws.Range('E31').Value = '  -2.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.85'
$ws.Range('E32').Value = '  -7.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.13'
$ws.Range('E33').Value = '  -7.42%  '
$ws.Range('E34').Value = '  -1.45%  '
$ws.Range('E35').Value = '  -7.12%  '
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('E37').Value = '  -7.39%  '
$ws.Range('E38').Value = '  -3.81%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.91'
$ws.Range('E39').Value = '  -3.35%  '
$ws.Range('E40').Value = '  +6.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0983'
$ws.Range('E41').Value = '  -6.17%  '
$ws.Range('E42').Value = '  -1.61%  '
$ws.Range('E43').Value = '  -1.10%  '
$ws.Range('E44').Value = '  -3.86%  '
$ws.Range('E45').Value = '  -5.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.50'
$ws.Range('E46').Value = '  -4.20%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.53'
$ws.Range('E47').Value = '  -5.76%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.87'
$ws.Range('E48').Value = '  -6.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.338.98'
$ws.Range('E49').Value = '  -5.49%  '
$ws.Range('E50').Value = '  -4.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.139.33'
$ws.Range('E51').Value = '  -4.61%  '
